# Apply the roster/odds reshuffle described by the diff.
# For each destination row, the F:V payload (teams, scores, odds, timestamps,
# url) is replaced by the F:V payload that currently lives in a different
# source row within the same match-day block. Columns A-E (index, country,
# league, season, kickoff date-time) stay attached to their row number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destRow -> srcRow (srcRow's current F:V values get copied into destRow)
$map = @{
    21 = 24; 22 = 23; 23 = 25; 24 = 22; 25 = 21;
    30 = 34; 31 = 33; 32 = 31; 33 = 30; 34 = 32;
    39 = 43; 40 = 39; 41 = 40; 42 = 41; 43 = 42;
    48 = 50; 49 = 48; 50 = 49; 51 = 52; 52 = 51;
    66 = 69; 67 = 70; 68 = 67; 69 = 66; 70 = 68;
    75 = 79; 76 = 78; 78 = 76; 79 = 75;
    87 = 88; 88 = 87;
    93 = 94; 94 = 95; 95 = 96; 96 = 93;
    102 = 103; 103 = 102; 104 = 105; 105 = 106; 106 = 104;
    137 = 138; 138 = 137;
}

# Snapshot the current F:V values of every row referenced (as source or
# destination) before writing anything back, since several blocks contain
# cyclic permutations (e.g. 21<-24, 24<-22, 22<-23, 23<-25, 25<-21).
$snapshot = @{}
foreach ($key in $map.Keys) {
    if (-not $snapshot.ContainsKey($key)) {
        $snapshot[$key] = $ws.Range("F$key`:V$key").Value()
    }
    $src = $map[$key]
    if (-not $snapshot.ContainsKey($src)) {
        $snapshot[$src] = $ws.Range("F$src`:V$src").Value()
    }
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $ws.Range("F$destRow`:V$destRow").Value = $snapshot[$srcRow]
}

# Append the five new matches (rows 139-143) completing match-day 17.
$ws.Range("A138:V138").Copy()
$ws.Range("A139:V143").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row = 139; A = 138; E = 45280.85416666666; F = "Stuttgart"; G = 3; H = "Augsburg"; I = 0;
       J = 1.45; K = "09/12/2023 20:03"; L = 1.38; M = "20/12/2023 20:29";
       N = 4.64; O = "09/12/2023 20:03"; P = 5.45; Q = "20/12/2023 20:29";
       R = 5.88; S = "09/12/2023 20:03"; T = 8.13; U = "20/12/2023 20:29";
       V = "https://www.betexplorer.com/football/germany/bundesliga/vfb-stuttgart-augsburg/2Lft1SyM/" },
    @{ Row = 140; A = 139; E = 45280.85416666666; F = "Bayer Leverkusen"; G = 4; H = "Bochum"; I = 0;
       J = 1.11; K = "09/12/2023 20:03"; L = 1.22; M = "20/12/2023 20:12";
       N = 9.17; O = "09/12/2023 20:03"; P = 7.42; Q = "20/12/2023 20:29";
       R = 14.27; S = "09/12/2023 20:03"; T = 12.08; U = "20/12/2023 20:28";
       V = "https://www.betexplorer.com/football/germany/bundesliga/bayer-leverkusen-bochum/KxFstRad/" },
    @{ Row = 141; A = 140; E = 45280.85416666666; F = "Eintracht Frankfurt"; G = 2; H = "B. Monchengladbach"; I = 1;
       J = 1.76; K = "09/12/2023 20:03"; L = 1.74; M = "20/12/2023 20:28";
       N = 3.93; O = "09/12/2023 20:03"; P = 4.21; Q = "20/12/2023 20:28";
       R = 4.06; S = "09/12/2023 20:03"; T = 4.63; U = "20/12/2023 20:29";
       V = "https://www.betexplorer.com/football/germany/bundesliga/eintracht-frankfurt-b-monchengladbach/Yu3P4Axc/" },
    @{ Row = 142; A = 141; E = 45280.85416666666; F = "Heidenheim"; G = 3; H = "Freiburg"; I = 2;
       J = 2.71; K = "09/12/2023 20:03"; L = 3.28; M = "20/12/2023 20:28";
       N = 3.42; O = "09/12/2023 20:03"; P = 3.43; Q = "20/12/2023 20:28";
       R = 2.49; S = "09/12/2023 20:03"; T = 2.32; U = "20/12/2023 20:28";
       V = "https://www.betexplorer.com/football/germany/bundesliga/heidenheim-freiburg/GObp0niS/" },
    @{ Row = 143; A = 142; E = 45280.85416666666; F = "Wolfsburg"; G = 1; H = "Bayern Munich"; I = 2;
       J = 5.98; K = "09/12/2023 20:03"; L = 6.44; M = "20/12/2023 20:29";
       N = 5.44; O = "09/12/2023 20:03"; P = 5.08; Q = "20/12/2023 20:28";
       R = 1.39; S = "09/12/2023 20:03"; T = 1.47; U = "20/12/2023 20:29";
       V = "https://www.betexplorer.com/football/germany/bundesliga/wolfsburg-bayern-munich/CW1T3Ui3/" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("A$r").Value = $nr.A
    $ws.Range("B$r").Value = "germany"
    $ws.Range("C$r").Value = "bundesliga"
    $ws.Range("D$r").Value = "2023-2024"
    $ws.Range("E$r").Value = $nr.E
    $ws.Range("F$r").Value = $nr.F
    $ws.Range("G$r").Value = $nr.G
    $ws.Range("H$r").Value = $nr.H
    $ws.Range("I$r").Value = $nr.I
    $ws.Range("J$r").Value = $nr.J
    $ws.Range("K$r").Value = $nr.K
    $ws.Range("L$r").Value = $nr.L
    $ws.Range("M$r").Value = $nr.M
    $ws.Range("N$r").Value = $nr.N
    $ws.Range("O$r").Value = $nr.O
    $ws.Range("P$r").Value = $nr.P
    $ws.Range("Q$r").Value = $nr.Q
    $ws.Range("R$r").Value = $nr.R
    $ws.Range("S$r").Value = $nr.S
    $ws.Range("T$r").Value = $nr.T
    $ws.Range("U$r").Value = $nr.U
    $ws.Range("V$r").Value = $nr.V
}
